$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (Price column D, Volume(1h) column E).
# Values are written with a leading apostrophe to force text type (matching
# the original inlineStr/text cells), then the style is reset to "Normal" so
# no quote-prefix style gets attached to the cell.
$updates = @{
    'D2' = '300.35'
    'E2' = '-1.15%'
    'D3' = '31.40'
    'E3' = '-1.99%'
    'D4' = '5.141'
    'E4' = '-2.68%'
    'D5' = '0.07337'
    'E5' = '-1.58%'
    'D6' = '1.812'
    'E6' = '23.31%'
    'D7' = '7.785'
    'E7' = '-0.81%'
    'D8' = '3.736'
    'E8' = '-0.81%'
    'D9' = '0.9244'
    'E9' = '0.98%'
    'D10' = '0.1688'
    'E10' = '-0.35%'
    'D11' = '0.07107'
    'E11' = '-7.76%'
    'D12' = '0.08113'
    'E12' = '1.90%'
    'D13' = '0.03021'
    'E13' = '-0.26%'
    'D14' = '0.09917'
    'E14' = '0.50%'
    'E15' = '-0.48%'
    'D16' = '0.006214'
    'E16' = '-1.40%'
    'E17' = '-0.76%'
    'D18' = '2.221'
    'E18' = '-0.27%'
    'D19' = '0.3226'
    'E19' = '-2.28%'
    'E20' = '-2.03%'
    'D21' = '4.551'
    'E21' = '0.62%'
    'D22' = '0.04651'
    'E22' = '1.75%'
    'D23' = '0.1581'
    'E23' = '-3.75%'
    'D24' = '0.001213'
    'D25' = '0.004749'
    'E25' = '7.57%'
    'D26' = '0.0001297'
    'E26' = '-7.50%'
    'D27' = '0.0001872'
    'E27' = '5.44%'
    'D39' = '0.01720'
    'E39' = '-1.57%'
    'E40' = '-0.14%'
    'D41' = '0.007083'
    'E41' = '-1.84%'
    'D42' = '0.1337'
    'E42' = '-0.07%'
    'E43' = '0.06%'
    'D44' = '0.01044'
    'E44' = '-23.04%'
    'D45' = '0.00006236'
    'E45' = '1.38%'
    'E46' = '-21.42%'
    'D47' = '1.920'
    'E47' = '2.54%'
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $updates[$cellRef]
    $range.Style = "Normal"
}
